$wb = $excel.ActiveWorkbook

# Sheet2 header row: rename the demo data titles
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "StringTitle"
$ws2.Range("B1").Value = "DateTitle"
$ws2.Range("C1").Value = "DoubleTitle"
$ws2.Range("D1").Value = "BigDecimalTitle"

# Update the saved selection on Sheet2 from D14 to D4
$ws2.Activate()
$null = $ws2.Range("D4").Select()
